$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manage_product")

# New "Details" description text for the Strawberry product (multi-line).
$details = @'
Details
Driscoll's is a renowned Californian berry producer. The family-owned company has been a staple in kitchens and in fruit bowls for over 100 years.
Amount Per Serving
Calories 32
'@
$details = $details.TrimEnd("`r", "`n")

# --- Header row: add "price" (D1) and "description" (E1) columns ---
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = "price"
$ws.Range("E1").Value = "description"

# --- Data row 2: add price value and description text ---
$ws.Range("C1").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = 300

$ws.Range("C1").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = $details
$ws.Range("E2").WrapText = $true
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("E2").VerticalAlignment = -4108

# Row height for row 2 to fit the wrapped description text.
$ws.Rows.Item(2).RowHeight = 68.25

# Column width for column E (description).
$ws.Columns.Item(5).ColumnWidth = 27.1
